$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.159.42"
$ws.Range("E2").Value = "  -4.34%  "
$ws.Range("D3").Value = "1.916.27"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "'245.02"
$ws.Range("E5").Value = "  -3.02%  "
$ws.Range("D6").Value = "'0.6986"
$ws.Range("E6").Value = "  -13.51%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("D8").Value = "'0.3217"
$ws.Range("E8").Value = "  -6.52%  "
$ws.Range("D9").Value = "'26.05"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'0.06788"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").Value = "'0.7886"
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("D12").Value = "'0.07931"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "1.920.64"
$ws.Range("E13").Value = "  -3.40%  "
$ws.Range("D14").Value = "'5.345"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "'93.55"
$ws.Range("E15").Value = "  -8.54%  "
$ws.Range("D16").Value = "'14.33"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "'258.89"
$ws.Range("E17").Value = "  -6.00%  "
$ws.Range("D18").Value = "30.169.26"
$ws.Range("E18").Value = "  -4.17%  "
$ws.Range("D19").Value = "'5.787"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "'0.000007787"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").Value = "2.175.90"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'6.784"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "'9.555"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'159.93"
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1309"
$ws.Range("E27").Value = "  -18.11%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'18.60"
$ws.Range("E28").Value = "  -5.93%  "
$ws.Range("D29").Value = "'2.201"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.544"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.336"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'4.379"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").Value = "'0.05016"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("D35").Value = "'1.181"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("D36").Value = "'0.7403"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "'0.01912"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").Value = "'2.783"
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("D40").Value = "'79.26"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").Value = "'6.479"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "'2.007"
$ws.Range("D43").Value = "'0.4367"
$ws.Range("E43").Value = "  -6.11%  "
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "'0.8331"
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("D46").Value = "'101.56"
$ws.Range("E46").Value = "  -4.07%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "'7.140"
$ws.Range("E48").Value = "  -4.92%  "
$ws.Range("D49").Value = "'35.63"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("D50").Value = "'0.05916"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "'1.466"
$ws.Range("E51").Value = "  +1.57%  "
